$wb = $excel.ActiveWorkbook

# --- 1) Goal sheet edits (these create new shared strings 619-624) ---
$goal = $wb.Worksheets.Item("Goal")

# Insert 2 new rows before row 32 so the old block (rows 32-41) shifts to rows 34-43
$goal.Rows("32:33").Insert()

# Row 28: "6." / "파일 첨부하는 기능."
$goal.Range("B28").Value = "6."
$goal.Range("C28").Value = "파일 첨부하는 기능."

# Row 30: date 42865 (2017-05-10) / "1." / "save 동작할때 item을 수정하도록 변경 제어 가능하도록.."
$goal.Range("A30").Value = 42865
$goal.Range("B30").Value = "1."
$goal.Range("C30").Value = "save 동작할때 item을 수정하도록 변경 제어 가능하도록.."

# Row 31: date 42867 (2017-05-12) / "2." / "RMA 정보 저장할때 날짜 저장될수 있도록."
$goal.Range("A31").Value = 42867
$goal.Range("B31").Value = "2."
$goal.Range("C31").Value = "RMA 정보 저장할때 날짜 저장될수 있도록."

# --- 2) RMA Table sheet edits (these create new shared strings 625-626) ---
$rma = $wb.Worksheets.Item("RMA Table")
$rma.Rows("4").Insert()
$rma.Range("A4").Value = "숫자형식"
$rma.Range("B4").Value = "DA 추가한 형식"

# Selection on RMA Table + make it the active (selected) tab
$rma.Activate()
$rma.Range("C1").Select()
